# Apply "Updated excel file with top entries" edit to metadata-template.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: was Name / # ['name'] -> now Format / # ['format'] ---
$ws.Range("A1").Value = "Format"
$ws.Range("B1").Value = "# ['format']"

# --- Row 2: Source / # ['source'] (unchanged content) ---
$ws.Range("A2").Value = "Source"
$ws.Range("B2").Value = "# ['source']"

# --- Row 3: new row -> Rows / #['rows'] ---
$ws.Range("A3").Value = "Rows"
$ws.Range("B3").Value = "#['rows']"

# --- Row 4 & 5: header row / template row, add Top / Moments columns ---
$ws.Range("A4").Value = "Column Name"
$ws.Range("B4").Value = "Column Type"
$ws.Range("C4").Value = "Missing"
$ws.Range("D4").Value = "No. of unique values"

$ws.Range("A5").Value = "# ['columns'][*]['name']"
$ws.Range("B5").Value = "# ['columns'][*]['type_pandas']"
$ws.Range("C5").Value = "# ['columns'][*]['missing']"
$ws.Range("D5").Value = "# ['columns'][*]['nunique']"

# Set column E (Top) pair before column F (Moments) pair so that the
# shared-string table is built in the same order as the target workbook.
$ws.Range("E4").Value = "Top"
$ws.Range("E5").Value = "#['columns'][*]['top']"

$ws.Range("F4").Value = "Moments"
$ws.Range("F5").Value = "#['columns'][*]['moments']"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 26.7109375
$ws.Columns.Item(2).ColumnWidth = 32
$ws.Columns.Item(3).ColumnWidth = 25
$ws.Columns.Item(4).ColumnWidth = 24.5703125
$ws.Columns.Item(5).ColumnWidth = 35.140625
$ws.Columns.Item(6).ColumnWidth = 24.42578125

# --- View / selection state ---
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("F6").Select()
